$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.929813666666668
$ws.Range("H2").Value = 26.789441
$ws.Range("I2").Value = 0.3579859341865942
$ws.Range("J2").Value = 0.3579859341865942
$ws.Range("M2").Value = 14.25737566666667
$ws.Range("N2").Value = 42.772127
$ws.Range("O2").Value = 0.2087950866344732
$ws.Range("P2").Value = 0.2087950866344732
$ws.Range("Q2").Value = 127.3157080790008
$ws.Range("R2").Value = 1145.841372711007
$ws.Range("S2").Value = 0.07474570414241274
$ws.Range("T2").Value = 0.07474570414241276

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.929813666666668
$ws.Range("H3").Value = 26.789441
$ws.Range("I3").Value = 0.3579859341865942
$ws.Range("J3").Value = 0.3579859341865942
$ws.Range("N3").Value = 87.128332
$ws.Range("O3").Value = 0.4253229592313036
$ws.Range("P3").Value = 0.4253229592313036
$ws.Range("Q3").Value = 259.3465899491569
$ws.Range("R3").Value = 2334.119309542412
$ws.Range("S3").Value = 0.1522596368914249
$ws.Range("T3").Value = 0.1522596368914249

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.929813666666668
$ws.Range("H4").Value = 26.789441
$ws.Range("I4").Value = 0.3579859341865942
$ws.Range("J4").Value = 0.3579859341865942
$ws.Range("M4").Value = 20.11084633333333
$ws.Range("N4").Value = 60.332539
$ws.Range("O4").Value = 0.2945174484164121
$ws.Range("P4").Value = 0.2945174484164122
$ws.Range("Q4").Value = 179.5861104356332
$ws.Range("R4").Value = 1616.274993920699
$ws.Range("S4").Value = 0.1054331039056014
$ws.Range("T4").Value = 0.1054331039056014

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.929813666666668
$ws.Range("H5").Value = 26.789441
$ws.Range("I5").Value = 0.3579859341865942
$ws.Range("J5").Value = 0.3579859341865942
$ws.Range("M5").Value = 4.873057999999999
$ws.Range("N5").Value = 14.619174
$ws.Range("O5").Value = 0.07136450571781097
$ws.Range("P5").Value = 0.07136450571781099
$ws.Range("Q5").Value = 43.51549992685933
$ws.Range("R5").Value = 391.6394993417341
$ws.Range("S5").Value = 0.02554748924715511
$ws.Range("T5").Value = 0.02554748924715511

# Row 6
$ws.Range("I6").Value = 0.03080543264277933
$ws.Range("J6").Value = 0.03080543264277933
$ws.Range("M6").Value = 14.25737566666667
$ws.Range("N6").Value = 42.772127
$ws.Range("O6").Value = 0.2087950866344732
$ws.Range("P6").Value = 0.2087950866344732
$ws.Range("Q6").Value = 10.955780926161
$ws.Range("R6").Value = 98.60202833544901
$ws.Range("S6").Value = 0.00643202297746154
$ws.Range("T6").Value = 0.00643202297746154

# Row 7
$ws.Range("I7").Value = 0.03080543264277933
$ws.Range("J7").Value = 0.03080543264277933
$ws.Range("N7").Value = 87.128332
$ws.Range("O7").Value = 0.4253229592313036
$ws.Range("P7").Value = 0.4253229592313036
$ws.Range("Q7").Value = 22.317312343476
$ws.Range("S7").Value = 0.0131022577720275
$ws.Range("T7").Value = 0.01310225777202751

# Row 8
$ws.Range("I8").Value = 0.03080543264277933
$ws.Range("J8").Value = 0.03080543264277933
$ws.Range("M8").Value = 20.11084633333333
$ws.Range("N8").Value = 60.332539
$ws.Range("O8").Value = 0.2945174484164121
$ws.Range("P8").Value = 0.2945174484164122
$ws.Range("Q8").Value = 15.453757537077
$ws.Range("R8").Value = 139.083817833693
$ws.Range("S8").Value = 0.009072737419315021
$ws.Range("T8").Value = 0.009072737419315023

# Row 9
$ws.Range("I9").Value = 0.03080543264277933
$ws.Range("J9").Value = 0.03080543264277933
$ws.Range("M9").Value = 4.873057999999999
$ws.Range("N9").Value = 14.619174
$ws.Range("O9").Value = 0.07136450571781097
$ws.Range("P9").Value = 0.07136450571781099
$ws.Range("Q9").Value = 3.744599085881999
$ws.Range("R9").Value = 33.701391772938
$ws.Range("S9").Value = 0.002198414473975267
$ws.Range("T9").Value = 0.002198414473975267

# Row 10
$ws.Range("G10").Value = 15.246351
$ws.Range("H10").Value = 45.739053
$ws.Range("I10").Value = 0.6112086331706265
$ws.Range("J10").Value = 0.6112086331706265
$ws.Range("M10").Value = 14.25737566666667
$ws.Range("N10").Value = 42.772127
$ws.Range("O10").Value = 0.2087950866344732
$ws.Range("P10").Value = 0.2087950866344732
$ws.Range("Q10").Value = 217.372953752859
$ws.Range("R10").Value = 1956.356583775731
$ws.Range("S10").Value = 0.1276173595145989
$ws.Range("T10").Value = 0.1276173595145989

# Row 11
$ws.Range("G11").Value = 15.246351
$ws.Range("H11").Value = 45.739053
$ws.Range("I11").Value = 0.6112086331706265
$ws.Range("J11").Value = 0.6112086331706265
$ws.Range("N11").Value = 87.128332
$ws.Range("O11").Value = 0.4253229592313036
$ws.Range("P11").Value = 0.4253229592313036
$ws.Range("Q11").Value = 442.7963772388439
$ws.Range("R11").Value = 3985.167395149596
$ws.Range("S11").Value = 0.2599610645678512
$ws.Range("T11").Value = 0.2599610645678512

# Row 12
$ws.Range("G12").Value = 15.246351
$ws.Range("H12").Value = 45.739053
$ws.Range("I12").Value = 0.6112086331706265
$ws.Range("J12").Value = 0.6112086331706265
$ws.Range("M12").Value = 20.11084633333333
$ws.Range("N12").Value = 60.332539
$ws.Range("O12").Value = 0.2945174484164121
$ws.Range("P12").Value = 0.2945174484164122
$ws.Range("Q12").Value = 306.617022105063
$ws.Range("R12").Value = 2759.553198945567
$ws.Range("S12").Value = 0.1800116070914957
$ws.Range("T12").Value = 0.1800116070914958

# Row 13
$ws.Range("G13").Value = 15.246351
$ws.Range("H13").Value = 45.739053
$ws.Range("I13").Value = 0.6112086331706265
$ws.Range("J13").Value = 0.6112086331706265
$ws.Range("M13").Value = 4.873057999999999
$ws.Range("N13").Value = 14.619174
$ws.Range("O13").Value = 0.07136450571781097
$ws.Range("P13").Value = 0.07136450571781099
$ws.Range("Q13").Value = 74.29635271135798
$ws.Range("R13").Value = 668.6671744022219
$ws.Range("S13").Value = 0.0436186019966806
$ws.Range("T13").Value = 0.04361860199668061

